$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("iOS App Store Review Guidelines")

$newRows = @(
    "Возможен бесплатный пробный период",
    "Возможно автоматическое продление подписки",
    "Подписки должны работать на всех устойствах пользователя",
    "При оформлении подписки основная функциональность остается",
    "Возможность обновления/отката",
    "При приобретении физичиских товаров/услуг возможен расчет через Apple Pay или кредитную карту",
    "Возможно хранение криптовалюты, облегчение транзакций или передачи криптовалюты",
    "Приложения не могут майнить криптовалюты",
    "Приложение должно включать функции, контент и пользовательский интерфейс, которые выводят его за пределы переупакованного веб-сайта",
    "Возможен вход через Apple ID",
    "Разрешение на доступ к ресурсам пользователя",
    "Возможность использовать приложение без входа учетную запись",
    "Собранные данные не могут быть перепрофилированы для других целей"
)

$tallRows = @(26, 29)

$startRow = 21
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $srcRow = $r - 1

    $ws.Range("A$srcRow" + ":B$srcRow").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("A$r").Formula = "=A$srcRow+1"
    $ws.Range("B$r").Value = $newRows[$i]

    if ($tallRows -contains $r) {
        $ws.Rows.Item($r).RowHeight = 30
    }
}

$excel.CutCopyMode = $false

$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("B20:B33").Select() | Out-Null
